# Update docs/epexspot_prices.xlsx: add a new date (08-sep / 2025-09-08)
# to the "Prix Spot" sheet, and two new trailing rows (2025-09-06,
# 2025-09-07) to the "Gaz" and "CO2" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Prix Spot" -> new column CI (08-sep) with 24 hourly values
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last existing column (CH) onto the new one
# (CI) so the header keeps the bold / bordered / centered style, then
# set the header text and the hourly values.
$wsSpot.Range("CH1:CH25").Copy()
$wsSpot.Range("CI1:CI25").PasteSpecial(-4122)

$wsSpot.Range("CI1").Value2 = "08-sep"

$spotValues = @(
    66.76000000000001,
    63.5,
    63.39,
    37.78,
    31.86,
    55.7,
    77.98,
    69.61,
    108.36,
    98.20999999999999,
    85.83,
    63,
    52.01,
    38.56,
    37.21,
    27.34,
    38.72,
    48.64,
    73.08,
    91.15000000000001,
    98.20999999999999,
    95.26000000000001,
    89.34999999999999,
    73.08
)

for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $wsSpot.Cells.Item($i + 2, 87).Value2 = $spotValues[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: "Gaz" -> two new rows (84, 85)
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Cells.Item(84, 1).Value2 = "'2025-09-06"
$wsGaz.Cells.Item(84, 1).Style = "Normal"
$wsGaz.Cells.Item(84, 2).Value2 = 31

$wsGaz.Cells.Item(85, 1).Value2 = "'2025-09-07"
$wsGaz.Cells.Item(85, 1).Style = "Normal"
$wsGaz.Cells.Item(85, 2).Value2 = 31

# ---------------------------------------------------------------------
# Sheet 3: "CO2" -> two new rows (84, 85)
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Cells.Item(84, 1).Value2 = "'2025-09-06"
$wsCo2.Cells.Item(84, 1).Style = "Normal"
$wsCo2.Cells.Item(84, 2).Value2 = 75.59

$wsCo2.Cells.Item(85, 1).Value2 = "'2025-09-07"
$wsCo2.Cells.Item(85, 1).Style = "Normal"
$wsCo2.Cells.Item(85, 2).Value2 = 75.59
